# Weekly price update: insert a new week's row of data at row 10,
# pushing the existing rows 10-13 down to rows 11-14.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 10 (shifts rows 10:13 -> 11:14, copies formatting down)
$ws.Rows("10:10").Insert()

# Populate the newly inserted row 10 with the new week's data
$ws.Cells.Item(10, 1).Value = 7
$ws.Cells.Item(10, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(10, 3).Value = "Ñuble"
$ws.Cells.Item(10, 4).Value = 44524
$ws.Cells.Item(10, 5).Value = 16
$ws.Cells.Item(10, 6).Value = 300000000
$ws.Cells.Item(10, 7).Value = "Espárragos"
$ws.Cells.Item(10, 8).Value = "Sin especificar"
$ws.Cells.Item(10, 9).Value = "Primera"
$ws.Cells.Item(10, 10).Value = 400
$ws.Cells.Item(10, 11).Value = 800
$ws.Cells.Item(10, 12).Value = 900
$ws.Cells.Item(10, 13).Value = 850
$ws.Cells.Item(10, 14).Value = "$/kilo"
$ws.Cells.Item(10, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(10, 16).Value = 850
$ws.Cells.Item(10, 17).Value = 1
$ws.Cells.Item(10, 18).Value = "Hortaliza"
